$wb = $excel.ActiveWorkbook

# --- 1. Insert a new worksheet "2022-Q3" before the existing "2022-Q1" sheet ---
$q1 = $wb.Worksheets.Item("2022-Q1")
$q3 = $wb.Worksheets.Add($q1)
$q3.Name = "2022-Q3"

# Copy the header-row + first-column formatting from "2022-Q1" so the new
# sheet matches its look (bordered/bold/centered header cells, etc.)
$q1.Range("B1:H1").Copy($q3.Range("B1:H1"))
$q1.Range("A2").Copy($q3.Range("A2:A4"))

# Header row text
$q3.Cells.Item(1,2).Value = "基金代码"
$q3.Cells.Item(1,3).Value = "基金名称"
$q3.Cells.Item(1,4).Value = "基金规模"
$q3.Cells.Item(1,5).Value = "股票总仓位"
$q3.Cells.Item(1,6).Value = "仓位占比"
$q3.Cells.Item(1,7).Value = "持有市值(亿元)"
$q3.Cells.Item(1,8).Value = "仓位排名"

# Row 2 (fund code / size / position values are kept as TEXT, matching the
# source sheet's formatting - a leading "'" forces text storage so values
# like "005457" / "0.0900" keep their leading/trailing zeros)
$q3.Cells.Item(2,1).Value = 0
$q3.Cells.Item(2,2).Value = "'005457"
$q3.Cells.Item(2,3).Value = "景顺长城量化小盘股票"
$q3.Cells.Item(2,4).Value = "'6.57"
$q3.Cells.Item(2,5).Value = "'93.58"
$q3.Cells.Item(2,6).Value = "'1.37"
$q3.Cells.Item(2,7).Value = "'0.0900"
$q3.Cells.Item(2,8).Value = 10

# Row 3
$q3.Cells.Item(3,1).Value = 1
$q3.Cells.Item(3,2).Value = "'015496"
$q3.Cells.Item(3,3).Value = "景顺中证1000指数增强C"
$q3.Cells.Item(3,4).Value = "'1.83"
$q3.Cells.Item(3,5).Value = "'92.63"
$q3.Cells.Item(3,6).Value = "'1.34"
$q3.Cells.Item(3,7).Value = "'0.0245"
$q3.Cells.Item(3,8).Value = 10

# Row 4
$q3.Cells.Item(4,1).Value = 2
$q3.Cells.Item(4,2).Value = "'015495"
$q3.Cells.Item(4,3).Value = "景顺中证1000指数增强A"
$q3.Cells.Item(4,4).Value = "'0.69"
$q3.Cells.Item(4,5).Value = "'92.63"
$q3.Cells.Item(4,6).Value = "'1.34"
$q3.Cells.Item(4,7).Value = "'0.0092"
$q3.Cells.Item(4,8).Value = 10

# --- 2. Update the "总计" summary sheet: add a new row for 2022-Q3 ---
#        (push the existing rows down: 2022-Q1 -> row3, 2021-Q1 -> row4)
$total = $wb.Worksheets.Item("总计")
$total.Range("A2").Copy($total.Range("A4"))

$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(4,2).Value = "2021-Q1"
$total.Cells.Item(4,3).Value = 1
$total.Cells.Item(4,4).Value = 0.14

$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(3,2).Value = "2022-Q1"
$total.Cells.Item(3,3).Value = 2
$total.Cells.Item(3,4).Value = 0.05

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q3"
$total.Cells.Item(2,3).Value = 3
$total.Cells.Item(2,4).Value = 0.12
